$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Text updates (order matters: it controls the order new entries
#    are appended to the shared-strings table so the final table
#    matches the target file's index layout).
# ------------------------------------------------------------------
$ws.Range("A16").Value = "Smoke Alarm Activation Time"
$ws.Range("A17").Value = "Sprinkler Activation Time"
$ws.Range("F1").Value = "RP, 1824 only"
$ws.Range("J1").Value = "RP all"
$ws.Range("B1").Value = "Original from NUREG-1934"

# ------------------------------------------------------------------
# 2) Styling: give the J1:L1 merged header and the K2:L2 sub-header
#    cells a "0.00" number format on top of their existing centered
#    alignment; give J2 the plain "0.00" numeric format used by the
#    rest of the data columns.
# ------------------------------------------------------------------
$ws.Range("J1:L1").NumberFormat = "0.00"
$ws.Range("J1:L1").HorizontalAlignment = -4108
$ws.Range("J1:L1").VerticalAlignment = -4108

$ws.Range("J2").NumberFormat = "0.00"
$ws.Range("K2:L2").NumberFormat = "0.00"

$ws.Range("J14").NumberFormat = "0.00"
$ws.Range("L14").NumberFormat = "0.00"

# ------------------------------------------------------------------
# 3) Updated comparison numbers (Dunes 200 rows) and McT-algorithm
#    removal from row 15.
# ------------------------------------------------------------------
$ws.Range("J4").Value = 0.98
$ws.Range("K4").Formula = "=0.45/2"

$ws.Range("J5").Value = 1.1599999999999999
$ws.Range("K5").Formula = "=0.43/2"

$ws.Range("J15").Value = 1

# ------------------------------------------------------------------
# 4) Remove the now-unused placeholder cells in rows 7, 10 and 11.
# ------------------------------------------------------------------
$ws.Range("J7:L7").Clear()
$ws.Range("J10:K10").Clear()
$ws.Range("J11:K11").Clear()

# ------------------------------------------------------------------
# 5) Fill in the new ATF rows (16 existing row gains data, 17 is a
#    brand-new row).
# ------------------------------------------------------------------
$ws.Range("J16").Value = 1.05
$ws.Range("K16").Formula = "=0.98/2"
$ws.Range("L16").Formula = "=0.33/2"

$ws.Range("J17").Value = 0.84
$ws.Range("J17").NumberFormat = "0.00"
$ws.Range("K17").Formula = "=0.52/2"
$ws.Range("K17").NumberFormat = "0.00"
$ws.Range("L17").Formula = "=0.2/2"
$ws.Range("L17").NumberFormat = "0.00"

# ------------------------------------------------------------------
# 6) Column width/format default for J:L (stored width "11" in the
#    sheet's character-width units) and the new selection.
# ------------------------------------------------------------------
$ws.Columns("J:L").ColumnWidth = 10.17

$ws.Range("B6").Select() | Out-Null
